$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 '64.532.39'
Set-TextValue 2 5 '  +0.13%  '
Set-TextValue 3 4 '3.143.70'
Set-TextValue 3 5 '  -1.16%  '
Set-TextValue 5 4 '573.22'
Set-TextValue 5 5 '  +0.31%  '
Set-TextValue 6 4 '164.46'
Set-TextValue 6 5 '  -2.71%  '
Set-TextValue 7 5 '  -0.11%  '
Set-TextValue 8 4 '0.577'
Set-TextValue 8 5 '  -5.27%  '
Set-TextValue 9 4 '3.161.19'
Set-TextValue 9 5 '  -0.89%  '
Set-TextValue 10 4 '0.118'
Set-TextValue 10 5 '  -2.56%  '
Set-TextValue 11 4 '6.64'
Set-TextValue 11 5 '  -2.35%  '
Set-TextValue 12 4 '0.385'
Set-TextValue 12 5 '  -0.87%  '
Set-TextValue 13 4 '3.695.40'
Set-TextValue 13 5 '  -1.16%  '
Set-TextValue 14 5 '  -1.79%  '
Set-TextValue 15 4 '64.513.97'
Set-TextValue 15 5 '  +0.02%  '
Set-TextValue 16 4 '25.07'
Set-TextValue 16 5 '  -1.07%  '
Set-TextValue 17 4 '3.148.55'
Set-TextValue 17 5 '  -0.90%  '
Set-TextValue 18 5 '  -2.17%  '
Set-TextValue 19 4 '407.74'
Set-TextValue 19 5 '  -2.80%  '
Set-TextValue 20 4 '5.27'
Set-TextValue 20 5 '  -1.90%  '
Set-TextValue 21 4 '12.57'
Set-TextValue 21 5 '  -3.12%  '
Set-TextValue 22 4 '7.10'
Set-TextValue 22 5 '  -0.39%  '
Set-TextValue 23 5 '  +0.12%  '
Set-TextValue 24 4 '69.06'
Set-TextValue 24 5 '  -1.88%  '
Set-TextValue 25 4 '0.485'
Set-TextValue 25 5 '  -1.10%  '
Set-TextValue 26 5 '  -5.36%  '
Set-TextValue 27 4 '0.0000103'
Set-TextValue 27 5 '  -2.93%  '
Set-TextValue 28 4 '8.93'
Set-TextValue 28 5 '  +2.08%  '
Set-TextValue 29 4 '0.994'
Set-TextValue 29 5 '  -0.69%  '
Set-TextValue 30 5 '  +0.12%  '
Set-TextValue 31 4 '1.81'
Set-TextValue 31 5 '  -1.69%  '
Set-TextValue 32 4 '21.28'
Set-TextValue 32 5 '  -2.29%  '
Set-TextValue 33 4 '162.40'
Set-TextValue 33 5 '  +3.78%  '
Set-TextValue 34 4 '4.88'
Set-TextValue 34 5 '  -3.33%  '
Set-TextValue 35 4 '6.31'
Set-TextValue 35 5 '  -0.77%  '
Set-TextValue 36 4 '1.13'
Set-TextValue 36 5 '  +0.08%  '
Set-TextValue 37 5 '  -0.07%  '
Set-TextValue 38 4 '1.69'
Set-TextValue 38 5 '  -0.91%  '
Set-TextValue 39 4 '2.645.01'
Set-TextValue 39 5 '  -1.89%  '
Set-TextValue 40 4 '23.79'
Set-TextValue 40 5 '  -1.86%  '
Set-TextValue 41 4 '4.10'
Set-TextValue 41 5 '  -2.73%  '
Set-TextValue 42 4 '38.28'
Set-TextValue 42 5 '  -2.52%  '
Set-TextValue 43 4 '0.692'
Set-TextValue 43 5 '  -3.34%  '
Set-TextValue 44 4 '0.0614'
Set-TextValue 44 5 '  -1.38%  '
Set-TextValue 45 4 '5.40'
Set-TextValue 45 5 '  -3.56%  '
Set-TextValue 46 2 'Bittensor'
Set-TextValue 46 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 46 4 '289.56'
Set-TextValue 46 5 '  -0.82%  '
Set-TextValue 47 2 'InjectiveProtocol'
Set-TextValue 47 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 47 4 '21.34'
Set-TextValue 47 5 '  -0.50%  '
Set-TextValue 48 2 'VeChain'
Set-TextValue 48 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 48 4 '0.0255'
Set-TextValue 48 5 '  -3.38%  '
Set-TextValue 49 5 '  -0.14%  '
Set-TextValue 50 4 '0.0977'
Set-TextValue 50 5 '  -1.39%  '
Set-TextValue 51 2 'dogwifhat'
Set-TextValue 51 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 51 4 '1.91'
Set-TextValue 51 5 '  -4.47%  '
